$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-TextCell($addr, $value) {
    # Prefix with a single quote so Excel stores the exact text
    # instead of reinterpreting numeric-looking strings as numbers.
    $ws.Range($addr).Value = "'" + $value
}

Set-Cell "D2" "27.710.44"
Set-Cell "E2" "  -4.78%  "

Set-Cell "D3" "1.723.72"
Set-Cell "E3" "  -5.10%  "

Set-TextCell "D4" "1.003"
Set-Cell "E4" "  -0.22%  "

Set-TextCell "D5" "225.37"
Set-Cell "E5" "  -3.79%  "

Set-TextCell "D6" "0.5718"
Set-Cell "E6" "  -3.63%  "

Set-Cell "E7" "  -0.17%  "

Set-TextCell "D8" "0.2710"
Set-Cell "E8" "  -0.76%  "

Set-TextCell "D9" "22.79"
Set-Cell "E9" "  -0.35%  "

Set-TextCell "D10" "0.06573"
Set-Cell "E10" "  -3.40%  "

Set-TextCell "D11" "0.07534"
Set-Cell "E11" "  -0.26%  "

Set-Cell "D12" "1.740.12"
Set-Cell "E12" "  -4.77%  "

Set-TextCell "D13" "4.668"
Set-Cell "E13" "  +0.12%  "

Set-TextCell "D14" "0.5961"
Set-Cell "E14" "  -3.76%  "

Set-Cell "D15" "1.966.71"
Set-Cell "E15" "  -4.63%  "

Set-TextCell "D16" "74.14"
Set-Cell "E16" "  -1.37%  "

Set-TextCell "D17" "0.000008590"
Set-Cell "E17" "  -10.16%  "

Set-Cell "D18" "27.707.82"
Set-Cell "E18" "  -4.31%  "

Set-TextCell "D19" "5.272"
Set-Cell "E19" "  -2.94%  "

Set-Cell "E20" "  -0.35%  "

Set-TextCell "D21" "203.68"
Set-Cell "E21" "  -2.24%  "

Set-TextCell "D22" "11.18"
Set-Cell "E22" "  -2.02%  "

Set-TextCell "D23" "6.567"
Set-Cell "E23" "  -2.85%  "

Set-TextCell "D24" "1.005"
Set-Cell "E24" "  -0.18%  "

Set-TextCell "D25" "149.36"
Set-Cell "E25" "  -3.31%  "

Set-TextCell "D26" "8.000"
Set-Cell "E26" "  +2.33%  "

Set-TextCell "D27" "0.1220"
Set-Cell "E27" "  -3.92%  "

Set-TextCell "D28" "16.09"
Set-Cell "E28" "  -0.99%  "

Set-Cell "B29" "Toncoin"
Set-Cell "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D29" "1.373"
Set-Cell "E29" "  -2.83%  "

Set-Cell "B30" "Hedera"
Set-Cell "C30" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D30" "0.06127"
Set-Cell "E30" "  -3.09%  "

Set-TextCell "D31" "1.386"
Set-Cell "E31" "  -3.56%  "

Set-TextCell "D32" "3.713"
Set-Cell "E32" "  -0.48%  "

Set-TextCell "D33" "3.691"
Set-Cell "E33" "  -0.18%  "

Set-TextCell "D34" "1.661"
Set-Cell "E34" "  -2.01%  "

Set-TextCell "D35" "1.026"
Set-Cell "E35" "  -3.95%  "

Set-TextCell "D36" "0.6403"
Set-Cell "E36" "  +0.72%  "

Set-TextCell "D37" "2.429"
Set-Cell "E37" "  -4.04%  "

Set-TextCell "D38" "2.677"
Set-Cell "E38" "  -2.75%  "

Set-TextCell "D39" "0.01653"
Set-Cell "E39" "  -3.86%  "

Set-Cell "D40" "1.114.51"
Set-Cell "E40" "  -1.13%  "

Set-TextCell "D41" "6.136"
Set-Cell "E41" "  -3.76%  "

Set-TextCell "D42" "0.8694"
Set-Cell "E42" "  +0.08%  "

Set-Cell "E43" "  -0.11%  "

Set-TextCell "D44" "99.39"
Set-Cell "E44" "  -0.80%  "

Set-Cell "D45" "1.878.42"
Set-Cell "E45" "  -4.71%  "

Set-TextCell "D46" "58.87"
Set-Cell "E46" "  -3.08%  "

Set-TextCell "D47" "0.00000000110"
Set-Cell "E47" "  -3.20%  "

Set-TextCell "D48" "1.551"
Set-Cell "E48" "  -2.05%  "

Set-TextCell "D49" "8.232"
Set-Cell "E49" "  -0.58%  "

Set-TextCell "D50" "0.05367"
Set-Cell "E50" "  -2.13%  "

Set-TextCell "D51" "0.4410"
Set-Cell "E51" "  -2.81%  "
